# Timesheet.xlsx edit — "more words for NaNoWriMo"
#
# Fills in previously-blank daily "Total" readings, lowers the daily goal
# delta, extends the tracking table (+ its Totals Row) by four more days,
# and widens the trend chart's plotted range to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Daily goal delta: 500 -> 400 (C2 feeds every Goal/D-column formula)
# ---------------------------------------------------------------------
$ws.Range("C2").Value = 400

# ---------------------------------------------------------------------
# 2. Backfill the "Total" column for rows that were logged later
# ---------------------------------------------------------------------
$totals = @{
    12 = 46818
    13 = 46818
    14 = 46818
    15 = 47253
    16 = 47502
    17 = 48843
    18 = 49118
    19 = 49118
    20 = 49439
    21 = 49439
    22 = 49439
    23 = 49439
    24 = 50028
    25 = 50306
    26 = 50750
}
foreach ($r in $totals.Keys) {
    $ws.Range("B$r").Value = $totals[$r]
}

# ---------------------------------------------------------------------
# 3. Append four more tracked days (rows 29-32), matching the format
#    and formulas already used by row 28 (copy format, then set
#    formula/value so no new cell styles get minted)
# ---------------------------------------------------------------------
$newDates = @{
    29 = 42314
    30 = 42315
    31 = 42316
    32 = 42317
}

for ($r = 29; $r -le 32; $r++) {
    $prev = $r - 1

    $ws.Range("A$prev").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)

    $ws.Range("B$prev").Copy()
    $ws.Range("B$r").PasteSpecial(-4122)

    $ws.Range("C$prev").Copy()
    $ws.Range("C$r").PasteSpecial(-4122)

    $ws.Range("D$prev").Copy()
    $ws.Range("D$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

foreach ($r in $newDates.Keys) {
    $ws.Range("A$r").Value = $newDates[$r]
}

for ($r = 29; $r -le 32; $r++) {
    $prev = $r - 1
    $ws.Range("C$r").Formula = "=B$r-B$prev"
    $ws.Range("D$r").Formula = "=D$prev+`$C`$2"
}

# ---------------------------------------------------------------------
# 4. Grow Table1 to cover the new rows and turn its Totals Row on
#    (this is what pushes the table ref out to A1:D33)
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D32"))
$lo.ShowTotals = $true

# Match formatting of the new totals row (row 33) to the rest of the
# table and make sure it is left blank (no per-column total function).
$ws.Range("A28").Copy()
$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("B28").Copy()
$ws.Range("B33").PasteSpecial(-4122)
$ws.Range("C28").Copy()
$ws.Range("C33").PasteSpecial(-4122)
$ws.Range("D28").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A33").Value = ""
$ws.Range("B33").Value = ""
$ws.Range("C33").Value = ""
$ws.Range("D33").Value = ""

# ---------------------------------------------------------------------
# 5. Widen the trend chart's plotted ranges to follow the table
# ---------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$s1 = $chart.SeriesCollection().Item(1)
$s2 = $chart.SeriesCollection().Item(2)
$s1.Formula = "=SERIES(Sheet1!`$B`$1,Sheet1!`$A`$2:`$A`$83,Sheet1!`$B`$2:`$B`$83,1)"
$s2.Formula = "=SERIES(Sheet1!`$D`$1,,Sheet1!`$D`$2:`$D`$83,2)"

# ---------------------------------------------------------------------
# 6. Leave the selection where the user was last working
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 21
$win.ScrollColumn = 1
$ws.Range("B27").Select()

$excel.CalculateFull()
